$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 62-64: add "Bots" label (column B) and a cost value of 0 (column G).
# These new cells pick up the same style used by the rest of the table
# (style index 1), so copy formatting from the existing, already-styled
# A-column cell in each row before writing the values.
62..64 | ForEach-Object {
    $row = $_
    $ws.Range("A$row").Copy()
    $ws.Range("B$row").PasteSpecial(-4122)
    $ws.Range("G$row").PasteSpecial(-4122)
    $ws.Range("B$row").Value = "Bots"
    $ws.Range("G$row").Value = 0
}

# Rows 65-85: fill in the previously-missing "cost" values (column G) for
# champions whose base cost was left blank in the original dataset.
$costRows   = @(65, 66, 67, 68, 69, 70, 71, 72, 73, 74, 75, 76, 77, 78, 79, 80, 81, 82, 83, 84, 85)
$costValues = @( 1,  1,  1,  2,  2,  2,  2,  3,  3,  3,  3,  3,  4,  4,  4,  4,  4,  4,  4,  5,  5)
for ($i = 0; $i -lt $costRows.Length; $i++) {
    $ws.Range("G" + $costRows[$i]).Value = $costValues[$i]
}

# Row 86 ("Target Dummy"): same "Bots" label as rows 62-64, plus a cost of 0.
$ws.Range("A86").Copy()
$ws.Range("B86").PasteSpecial(-4122)
$ws.Range("B86").Value = "Bots"
$ws.Range("G86").Value = 0

$excel.CutCopyMode = 0

# Restore the scrolled/selected view state recorded at save time.
$ws.Range("D39").Select()
$excel.ActiveWindow.ScrollRow = 25
